# MAVRIC Apparel.xlsx update
# - Fix header typo "SHirt Size" -> "Shirt Size"
# - Remove the blank "Chris Wong" row (was row 19)
# - Fill in shirt-size / T-shirt / Sweatshirt answers for Jake, John, Kyle and Derick
#   (Kyle's T-shirt answer also changes from "-" to "N")

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the empty "Chris Wong" row entirely - everything below shifts up one row
$ws.Rows.Item(19).Delete()

# After the deletion the rows are:
#   18 Alex Schell   (still blank)
#   19 Jake          (blank -> fill in)
#   20 John          (blank -> fill in)
#   21 Kyle          (already has data, T-shirt col changes)
#   22 Morgan        (already has data, unchanged)
#   23 Derick        (blank -> fill in)

# Copy the existing data-row formatting down onto the newly-populated rows
# so the cells pick up the same style (s="1") used throughout the table.
$ws.Range("B2:D2").Copy()
$ws.Range("B19:D20").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("B19").Value = "M"
$ws.Range("C19").Value = "Y"
$ws.Range("D19").Value = "N"

$ws.Range("B20").Value = "L"
$ws.Range("C20").Value = "Y"
$ws.Range("D20").Value = "N"

$ws.Range("C21").Value = "N"

$ws.Range("B23").Value = "L"
$ws.Range("C23").Value = "Y"
$ws.Range("D23").Value = "N"

# Fix the header typo in B1
$ws.Range("B1").Value = "Shirt Size"

# Restore the selection the sheet was saved with
$ws.Range("D29").Select() | Out-Null
